$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Japanese translation for HinaHead: "ヒナメット1" -> "ヒナメット"
$ws.Range("M7").Value = "ヒナメット"

# Add a new translation row for the Gusao skin (row 9 was empty before)
$ws.Range("A9").Value = "Gusao"
$ws.Range("M9").Value = "ぐさお(本体)"

# Move the "shito / 使徒" pair down one row (row 10 -> row 11)
$shitoKey = $ws.Range("A10").Value2
$shitoJp = $ws.Range("M10").Value2
$ws.Range("A11").Value = $shitoKey
$ws.Range("M11").Value = $shitoJp
$ws.Range("A10").ClearContents()
$ws.Range("M10").ClearContents()

# Move the "honentomori / 骨森" pair down one row (row 13 -> row 14)
$honenKey = $ws.Range("A13").Value2
$honenJp = $ws.Range("M13").Value2
$ws.Range("A14").Value = $honenKey
$ws.Range("M14").Value = $honenJp
$ws.Range("A13").ClearContents()
$ws.Range("M13").ClearContents()

# Update the active selection to match the final cursor position
$ws.Range("M13").Select() | Out-Null
